$d = $word.ActiveDocument

# Each "<id>...</id>" marker in this document is (pointlessly) split across
# three runs with different formatting:
#   run 1: "<id>"     - Courier New, color 7f6000, 9pt  (tag markup style)
#   run 2: "p039r_N"  - Arial, color 000000, 11pt        (plain body style)
#   run 3: "</id>"    - Courier New, color 7f6000, 9pt  (tag markup style)
# The fix collapses each triple into a single run so the whole tag
# "<id>p039r_N</id>" carries just the Courier New / 7f6000 / 9pt markup
# formatting, matching how the other freshly-downloaded tc/tcn/tl markers
# in the doc are represented. The visible text is unchanged - only the run
# boundaries/formatting are.

$wdReplaceAll = 2
$wdFindContinue = 1

$ids = @("p039r_1", "p039r_2")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"

    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, `
                                $wdFindContinue, $false, $needle, $wdReplaceAll)
    if (-not $found) {
        Write-Output "warning: '$needle' not found"
    }
}
